$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = '더 강력해진 SSG.COM 삼성카드'
$ws.Range("B3").Value = 'SSG MONEY 최대 15% 적립 + 스마일클럽 월이용료 할인'
$ws.Range("C3").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000002385&recruitmentPath=SSG'
$ws.Range("D3").Value = 'SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2022-10-26'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '2025-10-25'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = '[''이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''SSG머니 최대 15% 적립'', ''스마일클럽 월 이용료 3,900원 결제일 할인'', ''최대 15% SSG머니 적립 중 5%는 SSG.COM에서 제공하는 서비스로 자세한 내용은 SSG.COM 확인바람.'', ''01. SSG MONEY 최대 15% 적립!'', ''       이용실적 관계없이 적립 : 1 % + 이용실적 충족 시 적립 4% = 적립 가능한 최대 SSG머니 5%'', ''       카드혜택 + 스마일클럽 혜택'', ''       이용실적 관계없이 적립 1% + 이용실적 충족 시 적립 9% + 쓱배송/쓱배송 traders/새벽배송 상품 구매 시 5% = 적립 가능한 최대 SSG머니 15%'', ''최대 15% SSG MONEY 적립'', ''SSG.COM 삼성카드 최대 10% 적립+스마일클럽 5% 적립'', ''SSG.COM 삼성카드 최대 10% 적립(1% 적립+9% 추가 적립)'', ''1% 적립(전월 실적 조건 및 적립한도 없음)'', ''9% 추가 적립'', ''삼성카드 할인이 적용된 일시불 및 할부 이용금액은 제외됩니다.'', ''스마일클럽 5% 적립'', ''쓱배송/쓱배송 traders/새벽배송 상품 결제건에 한해 혜택을 받을 수 있습니다.'', ''02. 스마일클럽 월이용료 매월 3,900원 결제일 할인'', ''15,000원만 담아도 쓱 새벽배송 무료배송(SSG가입시)'', ''1쓱 새벽배송 트레이더 - SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''스마일클럽 단독 혜택 - 스타벅스 상품 전용 딜'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''카드 할인 혜택 자세히보기 (레이어팝업 열기)'', ''2022.10.26부터, SSG.COM 삼성카드 발급 시 스마일클럽에 자동 가입됩니다.'', ''G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우 혜택을 받을 수 없습니다.'', ''SSG.COM 삼성카드로 스마일클럽 월 이용료(3,900원) 정기결제 시 혜택이 제공됩니다.(월 1회)'', ''SSG.COM을 통한 스마일클럽 가입건에 한해 혜택이 제공됩니다.'', ''결제금액이 할인금액보다 적을 경우, 결제금액만큼 할인이 적용됩니다.(결제금액이 없는 경우 할인 대상에서 제외)'', ''01. 5만원 이상 결제 시 사용 가능한 4만 5천원 할인쿠폰 제공'', ''45,000원 할인쿠폰 - SSGPAY 바로결제 이용 시 사용 가능'', ''쿠폰 발급기간 : 2023.01.01 ~ 2023.01.31'', ''쿠폰 사용기간 : 2023.01.01 ~ 2023.01.31'', ''SSGPAY 바로결제 이용 이력이 없는 회원'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 50,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제하셔야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 통합 회원 본인 명의의 SSG.COM 삼성카드로 결제 시 사용 가능하며, 다른 부정적인 방법으로 사용한 경우에는 주문이 취소될 수 있습니다.'', ''쿠폰은 결제 화면에서 자동으로 적용됩니다.'', ''결제 화면에서 쿠폰 변경을 원할 경우 ‘쿠폰선택’을 눌러주세요.'', ''02. 스마일클럽 고객이라면 9% 장바구니 쿠폰'', ''SSG.COM 삼성카드로 100,000원 이상 결제 시 사용 가능(아이디당 1회, 할인한도 20,000원)'', ''SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 100,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 결제 시 최적 할인에 의해 자동 설정되며, 직접 변경할 수 있 습니다.'', ''현금성 상품, 무형서비스 상품, 초특가 상품 등 일부 상품은 제외됩니다.(쿠폰 적용 제외 상품은 결제 화면에서 쿠폰이 보이지 않음)'', ''03. 스마일클럽 월 이용료 결제 카드로 SSG.COM 삼성카드 등록 시 SSG머니 3,900원 적립'', ''(1)행사기간 동안 SSG.COM 삼성카드를 통해 스마일클럽 자동 가입 시 SSG머니 3,900원 즉시 적립'', ''적립시점까지 스마일클럽 월 이용료 정기결제 수단에 SSG.COM 삼성카드를 등록해야 혜택을 받을 수 있습니다.'', ''본 상품 발급 및 SSGPAY앱 설치시 바로결제에 자동 등록 됩니다.'', ''SSG MONEY는 매월 1일~말일까지 매출전표가 접수된 금액에 대해 다음달 25일 SSG.COM 계정으로 적립됩니다.'', ''적립된 SSG MONEY는 SSGPAY 회원가입 후 조회 및 사용 가능합니다. 단, SSG.COM 에서는 SSG.COM만 가입해도 사용 가능합니다.'', ''연체이자율 : 회원별/이용상품별 정상이자율+3.0%p(최고 연 20.0%)'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!스마일클럽 가입하고 모든 혜택 누리세요'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우, 월 이용료 결제일할인 혜택을 받을 수 없습니다.'']'

$ws.Range("A4").Value = '최대 5만원 혜택'
$ws.Range("B4").Value = '+ 쓸 때마다 최대 12% 적립'
$ws.Range("C4").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000000858&siteNo=6005&recruitmentPath=L6007001&eventCode=HPG02'
$ws.Range("D4").Value = 'SSG.COM카드 Edition 2 이벤트 안내 페이지'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2022-07-08'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '2023-06-28'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = '[''이벤트/쿠폰 > SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''SSG.COM카드 Edition2는 SSG.COM에서 최대 5만원 혜택 + 쓸 때마다 최대 12% 적립'', ''이달의 혜택 01. SSG.COM카드 Edition2 첫결제 시 1만원 쿠폰 할인 바로보기'', ''혜택 01. SSG머니 최대 12% 적립 바로보기'', ''혜택 02. 매월 스마일클럽 가입비 3,900원 지원 바로보기'', ''이달의 혜택 01. SSG.COM에서 SSG.COM카드 Edition2 첫 결제 시 1만원 할인'', ''      직전 6개월간 (2022년 7월 1일 부터 2022년 12월 31일) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한함'', ''첫 결제 쿠폰 이용 방법'', ''② SSGPAY에 등록된 SSG.COM카드 Edition2로'', ''③ SSG.COM에서 기간 내 1만 1천원 이상 첫 결제 시 1만원 쿠폰 할인'', ''      본 혜택은 SSG.COM카드 및 SSG.COM카드 Edition2 로 SSG.COM(이마트몰, 신세계몰, 신세계백화점몰 등)에서 직전 6개월간 (2022년 7월 1일 부터 년 12월 31일 까지) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한해 제공됨'', "      본 혜택은 결제 시 [결제방법 > SSGPAY카드] 내 SSG.COM카드 Edition2 선택 시, ''카드할인 최적가'' 추천에 의해 할인 금액이 자동 적용됨. 단, 1만 1천원 이상 첫 결제 시 적용)", ''다운받은 쿠폰은 SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''쿠폰 사용하러 가기'', ''SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''* 쿠폰할인, 적립금 사용 등 선할인을 제외한 카드 최종 결제금액이 2만원 이상인 경우에 한해 혜택 적용'', ''혜택 01. 장 볼 때마다 SSG머니 최대 12% 적립'', ''SSG.COM에서 최대 12%'', ''       스마일클럽으로 5% 적립'', ''       쓱·새벽·트레이더스 구매 시 (구매 전 스마일클럽 적립 아이콘을 꼭 확인해주세요)'', ''       SSG.COM 카드 Edition2로 7% 적립'', ''어디서나 한도없이 0.5%'', ''SSG.COM카드 Edition2로 어디서나 한도없이 0.5% 적립 (SSG.COM 외 모든 가맹점)'', ''SSG Money 최대 12% 적립'', ''스마일클럽 5% 적립 + SSG.COM카드 Edition2 최대 7% 적립'', ''      스마일클럽 5% 적립은 쓱·새벽배송·트레이더스 이용 시에 한함'', ''      SSG.COM카드 Edition2 최대 7% 적립'', ''      SSG.COM에서 결제 시 7%(1만 쓱머니 한도), 그 외 가맹점 0.5% 적립(적립한도 제한 없음)'', ''      무이자 할부 및 현대카드에서 제공하는 다른 할인 서비스 이용 시 적립 제외'', ''스마일클럽 가입비 3,900원 매월 100% 지원'', ''      1. 스마일클럽 자동 가입에'', ''      2. SSGPAY 내 카드 자동 등록'', ''      3. 스마일클럽 정기결제수단 자동 등록 및 월 이용료 3,900원 지원까지! (단, 해당 카드를 월 정기결제 수단에 등록한 경우에 한함)'', ''SSG.COM카드 Edition2를 스마일클럽 월 정기결제 수단에 등록 및 전월 이용금액 30만원 이상 시 혜택 제공'', ''TIP. 스마일클럽 가입 시 SSG.COM 혜택'', ''(SSG 가입 시) 15,000원만 담아도 쓱 · 새벽배송 무료배송'', ''쓱 · 새벽배송 · 트레이더스 SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''스마일클럽 단독 혜택 스타벅스 상품 전용 딜'', ''스마일클럽 가입비 매월 3,900원 지원'', ''월 1회, 매달 스마일클럽 정기결제일에 혜택 제공'', ''SSG.COM카드 Edition2는 최초 발급 시, 스마일클럽 월 정기결제 수단에 자동 등록 됨'', ''전월 이용금액 30만원 미만 시, SSG.COM카드 Edition2로 스마일클럽 정기 결제 금액이 자동 결제됨'', ''스마일클럽 무료 이용 기간이라면 정기결제 금액 지원 대신 SSG머니 3,900원 제공'', ''      스타벅스 자동 충전, 생활요금(통신요금, 아파트관리비 등) 정기결제 신청 및 이체 시 최대 1만원 청구 할인'', ''      2. 스타벅스 자동 충전 또는 생활요금 정기결제 신청(각 항목당 할인한도 5천원, 최대 1만원 할인)'', ''      정기결제 신청 후 카드 결제일에 따라 매출 발생 다음 달 또는 다다음 달 청구 할인 혜택 적용'', ''      단, 청구 할인 제공 일정은 당사 또는 신청인 사정에 의해 상이할 수 있음'', ''      3. 쏘카 1만원 할인쿠폰'', ''      쿠폰은 등록일 포함 30일간 이용 가능'', ''실물 SSG.COM카드 Edition2 수령 전 SSGPAY로 결제 시 건당 100만원 이하 결제 가능 *단, 본인 확인(신분증 확인 및 1원 인증) 완료한 경우에 한하며, 건당 100만원 초과 시 실물카드 수령 후 결제 가능'', ''SSG머니 최대 적립 12%에서 5%는 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능함'', ''카드 이용대금 연체 시 약정금리 + 연체가산금리 3%의 연체이자율이 적용됩니다. (회원별, 이용 상품별 차등적용 / 법정 최고금리 20% 이내) 단, 연체 발생시점에 약정금리가 없는 경우 아래와 같이 적용'', ''일시불 : 거래 발생시점 기준 최소 기간 (2개월)의 유이자 할부 약정금리 + 연체가산금리 3%'', ''무이자할부 : 거래발생시점 기준 동일한 할부 계약 기간의 유이자할부 약정금리 + 연체가산금리 3%'']'

$ws.Range("A5").Value = '쓱배송의 세계로 초대합니다'
$ws.Range("B5").Value = '첫 구매 쿠폰 받고, 친구 초대하면 혜택이 두 배'
$ws.Range("C5").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003900'
$ws.Range("D5").Value = '1월 쓱배송 친구 초대 이벤트'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2023-01-19'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '2023-01-31'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = '[''이벤트/쿠폰 > 1월 쓱배송 친구 초대 이벤트'', ''스마일클럽'', ''1월 쓱배송 친구 초대 이벤트'', ''5천원 할인에 무료배송'', ''친구 초대 이벤트'', ''쓱배송이 처음이라면? 할인에 무료배송'', ''쓱배송 웰컴 쿠폰'', ''5,000원 장바구니 쿠폰 - 첫구매 전용'', ''무료배송 - 첫구매 전용'', ''이마트몰, 쓱배송/새벽배송 20,000원 이상 구매시 사용 가능'', ''첫구매 쿠폰이 모두 소진되었습니다. 더 좋은 혜택으로 찾아뵙겠습니다.'', ''발급 대상 : 2022년 1월 18일 이후 이마트몰, 트레이더스 쓱배송/점포 택배 및 새벽배송 구매 이력이 없는 고객 '', ''쿠폰 사용 조건 : 이마트 쓱배송/새벽배송 상품 2만원 이상 구매 시 사용 가능'', ''이마트몰 쓱배송, 트레이더스 쓱배송, 새벽배송 세가지 배송 모두 한 번도 해본 적 없는 친구에게 나의 초대장 번호를 공유해주세요'', ''응모 이후 초대받은 친구가 이벤트 기간 내 쓱배송 첫 구매를 완료해야만 추천인과 초대받은 사람 모두에게 SSG MONEY 5천원이 적립됩니다. (*첫 구매 완료 기준 - 첫 구매 배송 완료 시점)'', ''SSG MONEY 적립 이후 초대받은 친구가 첫 쓱배송 주문을 전체 반품할 경우, 해당 초대 건에 대하여 SSG MONEY가 회수될 수 있습니다.'', ''본인 스스로 초대한 경우에는 이벤트 참여에서 제외됩니다.'', "친구가 전달 준 초대장 번호를 입력하고 ''첫 쓱배송 응모하기'' 버튼을 눌러주세요 응모 이후, 반드시 이벤트 기간 내 첫 쓱배송 구매를 완료해야만 SSG MONEY가 정상 적립됩니다.", ''첫 쓱배송 응모하기'', ''이벤트 응모여부 확인하기'', ''STEP2 첫 쓱배송 배송 주문하기 : 쓱배송 or 쓱배송TRADERS or 새벽배송(* 주문시, 쓱배송 혹은 새벽배송 스티커를 확인해주세요)'', ''STEP3 SSG MONEY적립 : 이벤트 종료 후 2월 10일(금) 이내 SSG MONEY가 일괄 적립됩니다.'', ''응모 이후 이벤트 기간 내 "쓱배송" 첫 구매를 완료해야만 추천인과 초대받은 사람 모두에게 SSG MONEY 5천원이 적립됩니다. (*첫 구매 완료 기준 - 첫 구매 배송 완료 시점)'', ''쓱배송이 아닌 다른 배송 유형으로 구매했을 시, 참여 인정되지 않습니다. (반드시 쓱배송 혹은 새벽배송 스티커를 확인해주세요)'', ''SSG MONEY 적립 이후 첫 쓱배송 주문을 전체 반품할 경우, 해당 초대 건에 대하여 SSG MONEY가 회수될 수 있습니다.'', ''이벤트 기간 :  1월 19일(목) 00:00 ~ 1월 31일(화) 23:59'', ''최근 1년 간 이마트몰 쓱배송, 트레이더스 쓱배송, 새벽배송 모두 구매 이력이 없는 경우만 쓱배송 첫 구매자로 응모 가능합니다. (2022년 1월 18일 이후 주문 내역이 없는 자)'', ''이벤트 기간 내 등록할 수 있는 초대자는 한 명입니다.'', ''본인 스스로 초대/응모한 경우에는 이벤트 참여에서 제외됩니다.'', ''응모 이후  초대받은 사람이 이벤트 기간 내 쓱배송 첫 구매를 완료해야만 추천인과 초대받은 사람 모두에게 SSG MONEY 5천원이 적립됩니다. (*첫 구매 완료 기준 ? 첫 구매 배송 완료)'', ''SSG MONEY 적립 이후 첫 쓱배송 주문을 전체 반품할 경우, 해당 초대 건에 대하여 추천인과 초대받은 사람 모두 SSG MONEY가 회수될 수 있습니다.'', ''초대 받은 친구가 탈퇴 후 재가입하여 구매를 하더라도 첫 구매로 인정되지 않으며, SSG MONEY가 지급되지 않습니다.'', ''본 이벤트를 통해 지급된 SSG MONEY는 지급일로부터 30일 이후 자동 소멸됩니다.'', ''부정한 방법으로 이벤트에 참여한 것이 발견될 경우, 당첨이 취소될 수 있습니다.'', ''본 이벤트에 대한 세부 사항은 당사 사정에 따라 임의로 변경 혹은 조기 종료될 수 있습니다.'', ''신세계상품권을 SSG MONEY로 전환하고, SSG.COM에서 쇼핑하세요'']'

$ws.Range("A6").Value = '2023 대한민국 수산대전'
$ws.Range("B6").Value = '설 특별전'
$ws.Range("C6").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003717'
$ws.Range("D6").Value = '(1/5~25) 2023 대한민국 수산대전 - 설 특별전'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2023-01-05'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '2023-01-25'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = '[''이벤트/쿠폰 > (1/5~25) 2023 대한민국 수산대전 - 설 특별전'', ''스마일클럽'', ''20% 장바구니 쿠폰'', ''여러 개 담으시고 최대 2만원 할인 받으세요!'', ''업체택배 상품 : 본 페이지 20% 장바구니 쿠폰 사용 가능 ID당 차수별 1매 발급, 최대 2만원 할인'', ''쓱배송 상품 더 보기'']'

$ws.Range("A7").Value = '스마일클럽 1월의 가입 혜택'
$ws.Range("B7").Value = 'WELCOME 5천원 할인쿠폰 제공'
$ws.Range("C7").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003449'
$ws.Range("D7").Value = '[스마일클럽] 1월 매일매일 스마일'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2023-01-02'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '2023-01-31'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = '[''이벤트/쿠폰 > [스마일클럽] 1월 매일매일 스마일'', ''스마일클럽'', ''[스마일클럽] 1월 매일매일 스마일'', '' 지금 스마일클럽 가입하고 모든 혜택 받기 '', ''한달 내내 막강한 스마일클럽 혜택'', ''WELCOME 스마일클럽이 처음이라면 첫달 무료 + 5천원 할인 쿠폰 발급 (5,100원 이상 구매 시 사용 가능)'', ''본 쿠폰은 1/5(목)~1/31(화) 기간 내 SSG.COM에서 스마일클럽에 최초 신규 가입한 고객님에 한해 ID당 1회 발급됩니다.'', ''쿠폰은 가입 차주 금요일 이내 자동 지급되며, 앱푸시나 문자메시지 등을 통해 별도 안내 예정입니다. (단, 마케팅 정보 수신을 비동의 하신 경우 안내 발송이 제한됩니다.)'', ''본 쿠폰은 5,100원 이상 구매 시 5,000원 할인됩니다. (할인액 및 배송비 제외한 구매 금액 기준으로 쿠폰 적용됨)'', ''일부 상품 및 브랜드는 쿠폰 적용 제외될 수 있습니다.'', ''SSG.COM의 쿠폰은 결제 시 최적 할인에 의해 자동 설정되며, 직접 변경이 가능합니다.'', ''본 이벤트는 당사 사정으로 내용이 변경되거나 종료될 수 있습니다.'', ''SSG.COM 카드 Edition2 매월 스마일클럽 가입비 100% 지원'', ''장 보는 날에도 멤버십은 더 큰 혜택 쓱배송데이 최대 2만원 10% 할인 쿠폰'', ''2023년 검은 토끼의 해 맞이 건강 프로젝트 스마일클럽 전용 최대 15% 쿠폰으로 건강을 선물해요'', ''계속해서 이어질 스마일클럽 전용 혜택을 기대해주세요'']'

$ws.Range("A8").Value = '쓱스럽게 안녕'
$ws.Range("B8").Value = 'See you again, Say Goodbye'
$ws.Range("C8").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003967'
$ws.Range("D8").Value = '쓱스럽게 안녕'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2023-01-17'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '2023-01-31'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = '[''이벤트/쿠폰 > 쓱스럽게 안녕'', ''스마일클럽'']'

$ws.Range("A9").Value = '1월 맘키즈 플러스'
$ws.Range("B9").Value = '매월 쏟아지는 맘키즈 특가!'
$ws.Range("C9").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000001665'
$ws.Range("D9").Value = '이달의 맘키즈 PLUS'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2022-09-01'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = '2999-12-13'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = '[''이벤트/쿠폰 > 이달의 맘키즈 PLUS'', ''스마일클럽'', ''맘키즈 ~40% 쿠폰상품'', ''지금 할인 중!'', ''※ 할인 금액은 상품 별로 상이하니 각 상품페이지를 꼭 참조하세요'', ''맘키즈 에누리 쿠폰'', ''행사대상 상품별 최대 40% 할인'', ''맘키즈 상품 에누리 쿠폰'', ''맘키즈 클럽 회원이라면 로그인 후 각 상품 상세페이지에서도 쿠폰을 받으실 수 있습니다.'', ''잠깐! 쿠폰이 적용되지 않는다면?'', ''맘키즈 에누리 쿠폰으로 상품별 최대 40%'', ''[파스퇴르] 올곧게만든 위드맘 3단계 750g (NEO2 쓱배송, 그외지역 택배)'', ''[파스퇴르] 올곧게만든 위드맘 2단계 750g (NEO2 쓱배송, 그외지역 택배)'', ''크리스탈 레진아트 칼라세트(N2 쓱배송, 전국택배)'', ''뽀로로 코딩컴퓨터(N2 쓱배송, 전국택배)'', ''미미 어린이병원 (N2쓱배송, 전국택배)'']'

$ws.Range("A10:G14").EntireRow.Delete()